{"js": "// Update the Weekly Metrics table row for \"Veatch\":\n//   Actuals from Previous Week: 160 -> 240\n//   Variance for Last Week:     +30 -> +80\n//   Forecast Time for Upcoming Week: 180 -> 140\n//   Totals to Date:             290 -> 530\n// and relocate the \"_GoBack\" bookmark (Word's \"last edit position\" marker)\n// from the trailing empty paragraph after the table to right after the\n// newly edited \"530\" value, matching where the author's last edit landed.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Row 3 (0-based) is the \"Veatch\" row; columns 1-4 hold the numbers.\nconst actualsCell = table.getCell(3, 1);\nconst varianceCell = table.getCell(3, 2);\nconst forecastCell = table.getCell(3, 3);\nconst totalsCell = table.getCell(3, 4);\n\nactualsCell.value = \"240\";\nvarianceCell.value = \"+80\";\nforecastCell.value = \"140\";\ntotalsCell.value = \"530\";\nawait context.sync();\n\n// Move the \"_GoBack\" bookmark to the end of the \"Totals to Date\" cell's\n// text (right after \"530\"), removing it from its previous location.\nconst existingBookmark = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\nexistingBookmark.load(\"isNullObject\");\nawait context.sync();\nif (!existingBookmark.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n}\n\nconst endOfTotals = totalsCell.body.getRange(\"End\");\nendOfTotals.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Update the Weekly Metrics table row for \"Veatch\":\n#   Actuals from Previous Week: 160 -> 240\n#   Variance for Last Week:     +30 -> +80\n#   Forecast Time for Upcoming Week: 180 -> 140\n#   Totals to Date:             290 -> 530\n# and relocate the \"_GoBack\" bookmark (Word's \"last edit position\" marker)\n# from the trailing empty paragraph after the table to right after the\n# newly edited \"530\" value, matching where the author's last edit landed.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Row 4 (1-based) is the \"Veatch\" row; columns 2-5 hold the numbers.\n$t.Cell(4, 2).Range.Text = \"240\"\n$t.Cell(4, 3).Range.Text = \"+80\"\n$t.Cell(4, 4).Range.Text = \"140\"\n$t.Cell(4, 5).Range.Text = \"530\"\n\n# Figure out the character position right after \"530\" (i.e. right before\n# the cell's end-of-cell mark) using a freshly re-fetched cell/range so we\n# see the edit that was just made.\n$totalsCell = $d.Tables.Item(1).Cell(4, 5)\n$endOfText = $totalsCell.Range.Duplicate()\n$endOfText.MoveEnd(1, -1) | Out-Null    # wdCharacter = 1; trim the cell mark\n$insertPos = $endOfText.End\n\n# Temporarily insert a marker character right after \"530\" so that the\n# target position is no longer the very last position in the paragraph\n# (collapsed bookmarks placed exactly at a paragraph/cell's final position\n# get relocated to the paragraph start, so we dodge that edge case).\n$endOfText.InsertAfter(\"#\")\n\n# A range built from absolute document positions resolves correctly for\n# Bookmarks.Add even in that edge case, so use that to plant the bookmark\n# immediately before the marker character (i.e. right after \"530\").\n$bookmarkRange = $d.Range($insertPos, $insertPos)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n\n# Remove the temporary marker character.\n$markerRange = $d.Range($insertPos, $insertPos + 1)\n$markerRange.Delete()\n"}
